# Updating timings based upon method inlining.
# The "Objeck (JIT)" column (C) got faster after inlining was applied,
# so its four sample timings - and everything derived from them
# (the C6 average and the A8 relative-delta formula) - need refreshing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.47585300000000003
$ws.Range("C3").Value = 0.48208499999999999
$ws.Range("C4").Value = 0.47908600000000001
$ws.Range("C5").Value = 0.47655199999999998

# Force a full recalculation so the AVERAGE/ABS formulas (C6, D6, A8) and
# the chart that plots row 6 all pick up the new numbers.
[void]$excel.CalculateFullRebuild()

$co = $ws.ChartObjects(1)
$chart = $co.Chart
[void]$chart.Refresh()

# Leave the selection where it ended up after reviewing the refreshed chart.
[void]$ws.Range("I5").Select()
